# "download articles with pandoc title blocks"
#
# The article used to open with an italic title line ("From Union Square to
# Rome, \n Chapter 7 - Reporting ====...") followed by a bold byline
# paragraph ("By Dorothy Day"). The new pandoc-style title block drops the
# title line entirely and turns the byline into a plain (non-bold) line
# prefixed with "%", e.g. "% Dorothy Day".

$d = $word.ActiveDocument

# 1. Remove the whole first paragraph (title line + line break + chapter
#    heading), paragraph mark included, so the byline paragraph becomes the
#    new first paragraph.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Delete()

# 2. Turn "By Dorothy Day" into "% Dorothy Day" with no bold formatting.
$bylinePara = $d.Paragraphs(1)
$bylineRange = $bylinePara.Range

$start = $bylineRange.Start
$end = $bylineRange.End

# Drop the paragraph-mark character from the range so we don't touch the
# paragraph mark's run properties, only the visible text run.
$textRange = $d.Range($start, $end - 1)
$textRange.Delete()

$newRange = $d.Range($start, $start)
$newRange.Text = "% Dorothy Day"
